$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected (accidental-edit guard); unprotect for the edits
# and re-protect afterwards so the workbook ends up in the same guarded state.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer banner (row 38, col A).
$oldText = $ws.Range("A38").Value2
$newText = $oldText -replace "2021-05-06", "2021-05-07"
$ws.Range("A38").Value = $newText

# Updated Weight (col D) and Percent Change (col E) figures for each holding row.
$rows = @(
    @{ Row = 2;  D = 0.03604417191491061;   E = 0.003881987577639689 }
    @{ Row = 3;  D = 0.02039925099396663;   E = 0.003501945525291816 }
    @{ Row = 4;  D = 0.01917348863762279;   E = 0.003170454270432721 }
    @{ Row = 5;  D = 0.03767643763347404;   E = 0.008075842696629199 }
    @{ Row = 6;  D = 0.03420385337455149;   E = 0.0004001600640257674 }
    @{ Row = 7;  D = 0.01977303184399096;   E = 0.0003859886133359414 }
    @{ Row = 8;  D = 0.0371772940331278;    E = 0.007424692405600508 }
    @{ Row = 9;  D = 0.02035040513705201;   E = 0.004500450045004545 }
    @{ Row = 10; D = 0.02557792944893501;   E = 0.003182812810821556 }
    @{ Row = 11; D = 0.02395227827349544;   E = 0.006903876792352603 }
    @{ Row = 12; D = 0.05753355049678938;   E = 0.008028335301062395 }
    @{ Row = 13; D = 0.02497549721365475;   E = 0.005867253392005667 }
    @{ Row = 14; D = 0.02755771308700555;   E = 0.003179409538228928 }
    @{ Row = 15; D = 0.03378749332545533;   E = -0.009613762860516184 }
    @{ Row = 16; D = 0.01934209435947186;   E = 0.009951518244450064 }
    @{ Row = 17; D = 0.03077288985620827;   E = 0.01430224867724861 }
    @{ Row = 18; D = 0.04195655584561729;   E = 0.003917050691244262 }
    @{ Row = 19; D = 0.1258214322531487;    E = 0.003984063745019695 }
    @{ Row = 20; D = 0.009102323673836841;  E = 0.008099768578040401 }
    @{ Row = 21; D = 0.01539819846242422;   E = 0.002775657483866434 }
    @{ Row = 22; D = 0.01671917360535887;   E = 0.01158880438719034 }
    @{ Row = 23; D = 0.0156398328108487;    E = 0.003191489361702216 }
    @{ Row = 24; D = 0.02127282861622397;   E = 0.01398745243237665 }
    @{ Row = 25; D = 0.01242251903540598;   E = -0.006311745334796837 }
    @{ Row = 26; D = 0.04221197897239997;   E = 0.001627251030592358 }
    @{ Row = 27; D = 0.02385936938315576;   E = 0.00009809691975681822 }
    @{ Row = 28; D = 0.04551823291230806;   E = 0.00619047619047608 }
    @{ Row = 29; D = 0.05592687797200517;   E = 0.008024251069900235 }
    @{ Row = 30; D = 0.01317616990271775;   E = 0.01346153846153841 }
    @{ Row = 31; D = 0.02062338224366338;   E = -0.0003824091778201755 }
    @{ Row = 32; D = 0.01348349175247221;   E = 0.01320754716981143 }
    @{ Row = 33; D = 0.04181861718093445;   E = 0.001029336078229637 }
    @{ Row = 34; D = 0.01675163574776671;   E = 0.01518083048072616 }
    @{ Row = 35; D = 0.9999999999999999;    E = 0.004907366177980155 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# Restore sheet protection (objects/scenarios protected, row/column formatting allowed).
$ws.Protect($null, $true, $true, $true, $false, $true, $false, $false, $true, $true, $true, $true, $true, $true, $true, $true)
